# remove phone # from resume
#
# The header line reads:
#   PHONE (226) 989-7867  •  E-MAIL JUDE.FIORILLO@GMAIL.COM  •  EHJUDE.COM
# and should become:
#   E-MAIL JUDE.FIORILLO@GMAIL.COM  •  EHJUDE.COM
#
# Word also records the location of the most recent edit with the
# reserved "_GoBack" bookmark, so after deleting the phone number we
# relocate that bookmark to sit right before the (now first) "E-MAIL"
# text, removing it from its old spot later in the document.

$d = $word.ActiveDocument

# 1) Remove "PHONE (226) 989-7867  " (including the two trailing spaces)
$rngPhone = $d.Content
[void]$rngPhone.Find.Execute("PHONE (226) 989-7867  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngPhone.Delete()

# 2) Remove the bullet separator ("\u2022  ") that used to sit between the
#    phone number and the e-mail address, so the line now starts at E-MAIL.
$rngBullet = $d.Content
[void]$rngBullet.Find.Execute([char]0x2022 + "  ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngBullet.Delete()

# 3) Move the "_GoBack" bookmark from its old location onto the start of
#    the (now leading) "E-MAIL" text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rngEmail = $d.Content
[void]$rngEmail.Find.Execute("E-MAIL", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkPoint = $d.Range($rngEmail.Start, $rngEmail.Start)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkPoint)
